$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,16,17,21,22,23,25,36,37,41,42,43,45,56,57,61,62,63,65,76,77,81,82,83,84,85,86,95,96,100,101,102,103,104,105,114,115,119,120,121,122,123,124,133,134,138,139,140,141,142,143,152,153,157,158,159,160,161,162,171,172,176,177,178,180,191,192,196,197,198,200,211,212,216,217,218,220,231,232)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
